# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2410   (the "before" format version)
#   *_new -> *_FV2504   (the "after" format version)
# and turn the header range into a real Excel Table ("Table1") so the new
# column names are picked up as the table's column headers, plus freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Rename the header cells (A1:U1). This rewrites the underlying shared
#    string text used by the header row - diff (K1) is left untouched.
# ---------------------------------------------------------------------
$headerNames = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headerNames.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerNames[$i]
}

# ---------------------------------------------------------------------
# 2. Stash the existing header formatting (bold font / grey fill / border /
#    centered+wrapped) on a scratch row, clear direct formatting from the
#    header row so turning it into a Table does not bake that formatting
#    into a header-row dxf, create the Table, then restore the formatting.
# ---------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A200:U200")

$headerRange.Copy($scratchRange)
$headerRange.ClearFormats()

$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U79"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)  # xlPasteFormats
$scratchRange.Clear()

# ---------------------------------------------------------------------
# 3. Freeze the header row (pane split below row 1).
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

Write-Output "done"
